# add react + express
# This script reproduces the edits made to data/transaction.xlsx:
#  - adds new transaction rows (12-18, including a SELL trade) to "transaction"
#  - adds network_fee / wallet_received columns (H..L) to "transaction"
#  - adds new deposit rows (4-7) to "deposit"
#  - updates selections on both sheets
#
# NOTE: cell writes below are intentionally ordered to match the original
# authoring session so that new shared-string entries land at the same
# index as in the source workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 2: "deposit" (populated first so that sheet 1 - "transaction" - ends
# up as the active/selected tab, matching the source workbook)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("deposit")

$depositRows = @(
    @{ Row=4; Date=44344; From="Citibank"; To="Coinbase "; Amount=500 },
    @{ Row=5; Date=44345; From="Citibank"; To="Coinbase "; Amount=210 },
    @{ Row=6; Date=44349; From="Citibank"; To="Webull";    Amount=113 },
    @{ Row=7; Date=44353; From="Citibank"; To="Coinbase "; Amount=500 }
)

foreach ($item in $depositRows) {
    $r = $item.Row
    $ws2.Cells.Item($r, 1).Value = $item.Date
    $ws2.Cells.Item($r, 2).Value = $item.From
    $ws2.Cells.Item($r, 3).Value = $item.To
    $ws2.Cells.Item($r, 4).Value = $item.Amount
}

# Update selection on sheet 2
$ws2.Range("E7").Select()

# ---------------------------------------------------------------------------
# Sheet 1: "transaction"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("transaction")

# --- New transaction rows 12-16 -------------------------------------------
$newRows = @(
    @{ Row=12; Date=44344; Symbol="ADA";  Name="Cardano";             Side="BUY";  Qty=1.5422400000000001; Price=161 },
    @{ Row=13; Date=44345; Symbol="ADA";  Name="Cardano";             Side="BUY";  Qty=1.365;               Price=153.9 },
    @{ Row=14; Date=44345; Symbol="ADA";  Name="Cardano";             Side="BUY";  Qty=1.4;                 Price=178.2 },
    @{ Row=15; Date=44349; Symbol="ARKK"; Name="ARK Innovation ETF";  Side="SELL"; Qty=111;                 Price=1 },
    @{ Row=16; Date=44350; Symbol="RIOT"; Name="Riot Blockchain";     Side="BUY";  Qty=28.5;                Price=5 }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws1.Cells.Item($r, 1).Value = $item.Date
    $ws1.Cells.Item($r, 2).Value = $item.Symbol
    $ws1.Cells.Item($r, 3).Value = $item.Name
    $ws1.Cells.Item($r, 4).Value = $item.Side
    $ws1.Cells.Item($r, 5).Value = $item.Qty
    $ws1.Cells.Item($r, 6).Value = $item.Price
    $ws1.Cells.Item($r, 7).Formula = "=E$r*F$r"
}

# Row 17 - BTC buy with real network-fee figures
$ws1.Cells.Item(17, 1).Value = 44353
$ws1.Cells.Item(17, 2).Value = "BTC"
$ws1.Cells.Item(17, 3).Value = "Bitcoin"
$ws1.Cells.Item(17, 4).Value = "BUY"
$ws1.Cells.Item(17, 5).Value = 36056.870000000003
$ws1.Cells.Item(17, 6).Value = 0.0069095399999999996
$ws1.Cells.Item(17, 7).Formula = "=E17*F17"

# Row 18 - BTC buy with real network-fee figures
$ws1.Cells.Item(18, 1).Value = 44353
$ws1.Cells.Item(18, 2).Value = "BTC"
$ws1.Cells.Item(18, 3).Value = "Bitcoin"
$ws1.Cells.Item(18, 4).Value = "BUY"
$ws1.Cells.Item(18, 5).Value = 36147.65
$ws1.Cells.Item(18, 6).Value = 0.0068898300000000004
$ws1.Cells.Item(18, 7).Formula = "=E18*F18"

# --- New column headers (H, K, I, J, L order matches original authoring) --
$ws1.Cells.Item(1, 8).Value  = "network_fee"
$ws1.Cells.Item(1, 11).Value = "network_fee_transfer_to_wallet"
$ws1.Cells.Item(1, 9).Value  = "network_fee_pct"
$ws1.Cells.Item(1, 10).Value = "wallet_received"
$ws1.Cells.Item(1, 12).Value = "network_fee_transfer_to_wallet_pct"

# New column widths (H..L)
$ws1.Columns.Item(8).ColumnWidth  = 13
$ws1.Columns.Item(9).ColumnWidth  = 15.33
$ws1.Columns.Item(10).ColumnWidth = 16
$ws1.Columns.Item(11).ColumnWidth = 31.5
$ws1.Columns.Item(12).ColumnWidth = 37.17

# Fill in H:L for rows 2-16 with the default 0 / 0 / #DIV0! pattern
for ($r = 2; $r -le 16; $r++) {
    $ws1.Cells.Item($r, 8).Value  = 0
    $ws1.Cells.Item($r, 9).Formula = "=(H$r/F$r)*100"
    $ws1.Cells.Item($r, 10).Value = 0
    $ws1.Cells.Item($r, 11).Value = 0
    $ws1.Cells.Item($r, 12).Formula = "=(K$r/J$r)*100"
}

# Row 17 network-fee figures
$ws1.Cells.Item(17, 8).Formula = "=F17-0.00689577"
$ws1.Cells.Item(17, 9).Formula = "=(H17/F17)*100"
$ws1.Cells.Item(17, 10).Value = 0.0068957699999999999
$ws1.Cells.Item(17, 11).Formula = "=F17-J17"
$ws1.Cells.Item(17, 12).Formula = "=(K17/J17)*100"

# Row 18 network-fee figures
$ws1.Cells.Item(18, 8).Value = 0.000024070000000000002
$ws1.Cells.Item(18, 9).Formula = "=(H18/F18)*100"
$ws1.Cells.Item(18, 10).Value = 0.0068790500000000003
$ws1.Cells.Item(18, 11).Formula = "=F18-J18"
$ws1.Cells.Item(18, 12).Formula = "=(K18/J18)*100"

# Update selection on sheet 1 (kept as the final/active sheet & selection)
$ws1.Range("C21").Select()
